# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated K (column G) values for the corresponding rows
$kValues = @{
    2  = 2
    3  = 1
    4  = 3
    5  = 2
    6  = 3
    7  = 1
    8  = 2
    9  = 0
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 2
    17 = 0
    19 = 1
    20 = 1
    21 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
